$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed NATMI/TPM run now reports a new "ECs" target cluster that
# sorts before the pre-existing "FAPs" cluster, so the original FAPs data
# row is pushed down to row 3 (with recomputed values), and row 2 now
# holds the new ECs row.

# Row 2: Resolving-Mac / Ccl22 -> Ccr4 / ECs (new target cluster, replaces old FAPs row in place)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.945012333333333
$ws.Range("H2").Value = 5.835037
$ws.Range("M2").Value = 0.04581866666666667
$ws.Range("N2").Value = 0.137456
$ws.Range("O2").Value = 0.4376255591461182
$ws.Range("P2").Value = 0.4376255591461182
$ws.Range("Q2").Value = 0.08911787176355555
$ws.Range("R2").Value = 0.802060845872
$ws.Range("S2").Value = 0.4376255591461182
$ws.Range("T2").Value = 0.4376255591461182

# Row 3 (new): Resolving-Mac / Ccl22 -> Ccr4 / FAPs (recomputed values for the pre-existing cluster)
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Ccl22"
$ws.Range("C3").Value = "Ccr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.945012333333333
$ws.Range("H3").Value = 5.835037
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05887966666666666
$ws.Range("N3").Value = 0.176639
$ws.Range("O3").Value = 0.5623744408538818
$ws.Range("P3").Value = 0.5623744408538818
$ws.Range("Q3").Value = 0.1145216778492222
$ws.Range("R3").Value = 1.030695100643
$ws.Range("S3").Value = 0.5623744408538818
$ws.Range("T3").Value = 0.5623744408538818

$wb.Save()
